$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the date column from an existing row, then set new value
$ws.Range("A2").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 42604.424884259257

$ws.Range("B7").Value = "Gilead Sciences, Inc."
$ws.Range("C7").Value = "GILD"

$ws.Range("D7").Value = 80.56
$ws.Range("E7").Value = 81
$ws.Range("F7").Value = 0.51
$ws.Range("G7").Value = 80.59
